$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the quantity (D) values and dates (A) for the three data rows,
# also fix up the SKU/unit text in row 4 per the web-service refresh.
$ws.Range("A2").Value = 42529
$ws.Range("A3").Value = 42530
$ws.Range("A4").Value = 42531

$ws.Range("D2").Value = 99
$ws.Range("D3").Value = 99
$ws.Range("D4").Value = 99

$ws.Range("C4").Value = "B3"
$ws.Range("E3").Value = "kg/ltr"
$ws.Range("E4").Value = "kg/ltr"

# Custom date display format for the Month column.
$ws.Range("A2:A4").NumberFormat = "dd\-mmm\-yyyy"

# Column widths re-fit for the refreshed data (best-fit, as Excel recalculates
# on edit).
$ws.Columns.Item(2).ColumnWidth = 15.59
$ws.Columns.Item(3).ColumnWidth = 16.02
$ws.Columns.Item(4).ColumnWidth = 18.45
$ws.Columns.Item(5).ColumnWidth = 14.74

# Move the active selection to D2.
[void]$ws.Range("D2").Select()
